$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 289, shifting rows 289-307 down to 290-308
$ws.Rows.Item(289).Insert()

# Populate the newly inserted row 289 with the new weekly data point.
# Columns A,B,C,E,F,G,H,I,N,Q,R are constant across this block of records.
$ws.Range("A289").Value = 10
$ws.Range("B289").Value = 'Vega Modelo de Temuco'
$ws.Range("C289").Value = 'La Araucanía'
$ws.Range("D289").Value = 44931
$ws.Range("E289").Value = 9
$ws.Range("F289").Value = 100112052
$ws.Range("G289").Value = 'Albahaca'
$ws.Range("H289").Value = 'Sin especificar'
$ws.Range("I289").Value = 'Primera'
$ws.Range("J289").Value = 300
$ws.Range("K289").Value = 6000
$ws.Range("L289").Value = 6000
$ws.Range("M289").Value = 6000
$ws.Range("N289").Value = '$/paquete'
$ws.Range("O289").Value = 'Región del Maule'
$ws.Range("P289").Value = 6000
$ws.Range("Q289").Value = 1
$ws.Range("R289").Value = 'Hortaliza'
